$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 1; this shifts all existing data
# (and its formatting) down by one row.
$ws.Rows("1").Insert()

# The insert operation carried the old header formatting down onto the
# data that moved to row 2. Move that formatting back up onto the new
# row 1 (which will hold the numeric header row) and strip it from row 2
# (which now holds the old, unstyled header text).
$ws.Range("A2:H2").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2:H2").ClearFormats()

# Populate the new row 1 with the numeric header values 0-7.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7

# The old "thread_size"/"material_surface" header text (originally G1/H1)
# is dropped rather than carried down to G2/H2.
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
